$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data.
# Force each touched cell to Text format before writing so that
# numeric-looking strings (e.g. "586.85") are preserved verbatim
# as text instead of being auto-converted to numbers, then restore
# the default "Normal" style so no stray formatting is introduced.

$cell = $ws.Range('D2')
$cell.NumberFormat = '@'
$cell.Value = '67.658.92'
$cell.Style = 'Normal'

$cell = $ws.Range('E2')
$cell.NumberFormat = '@'
$cell.Value = '  +0.90%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D3')
$cell.NumberFormat = '@'
$cell.Value = '2.492.21'
$cell.Style = 'Normal'

$cell = $ws.Range('E3')
$cell.NumberFormat = '@'
$cell.Value = '  +0.78%  '
$cell.Style = 'Normal'

$cell = $ws.Range('E4')
$cell.NumberFormat = '@'
$cell.Value = '  +0.03%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D5')
$cell.NumberFormat = '@'
$cell.Value = '586.85'
$cell.Style = 'Normal'

$cell = $ws.Range('E5')
$cell.NumberFormat = '@'
$cell.Value = '  +0.71%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D6')
$cell.NumberFormat = '@'
$cell.Value = '176.48'
$cell.Style = 'Normal'

$cell = $ws.Range('E6')
$cell.NumberFormat = '@'
$cell.Value = '  +4.36%  '
$cell.Style = 'Normal'

$cell = $ws.Range('E7')
$cell.NumberFormat = '@'
$cell.Value = '  -0.02%  '
$cell.Style = 'Normal'

$cell = $ws.Range('E8')
$cell.NumberFormat = '@'
$cell.Value = '  +0.27%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D9')
$cell.NumberFormat = '@'
$cell.Value = '0.139'
$cell.Style = 'Normal'

$cell = $ws.Range('E9')
$cell.NumberFormat = '@'
$cell.Value = '  +3.76%  '
$cell.Style = 'Normal'

$cell = $ws.Range('E10')
$cell.NumberFormat = '@'
$cell.Value = '  +0.25%  '
$cell.Style = 'Normal'

$cell = $ws.Range('E11')
$cell.NumberFormat = '@'
$cell.Value = '  +2.35%  '
$cell.Style = 'Normal'

$cell = $ws.Range('E12')
$cell.NumberFormat = '@'
$cell.Value = '  +0.14%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D13')
$cell.NumberFormat = '@'
$cell.Value = '2.949.88'
$cell.Style = 'Normal'

$cell = $ws.Range('E13')
$cell.NumberFormat = '@'
$cell.Value = '  +0.86%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D14')
$cell.NumberFormat = '@'
$cell.Value = '25.72'
$cell.Style = 'Normal'

$cell = $ws.Range('E14')
$cell.NumberFormat = '@'
$cell.Value = '  +0.83%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D15')
$cell.NumberFormat = '@'
$cell.Value = '67.494.88'
$cell.Style = 'Normal'

$cell = $ws.Range('E15')
$cell.NumberFormat = '@'
$cell.Value = '  +1.00%  '
$cell.Style = 'Normal'

$cell = $ws.Range('E16')
$cell.NumberFormat = '@'
$cell.Value = '  +1.47%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D17')
$cell.NumberFormat = '@'
$cell.Value = '2.485.44'
$cell.Style = 'Normal'

$cell = $ws.Range('E17')
$cell.NumberFormat = '@'
$cell.Value = '  +2.38%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D18')
$cell.NumberFormat = '@'
$cell.Value = '11.07'
$cell.Style = 'Normal'

$cell = $ws.Range('E18')
$cell.NumberFormat = '@'
$cell.Value = '  +0.21%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D19')
$cell.NumberFormat = '@'
$cell.Value = '7.46'
$cell.Style = 'Normal'

$cell = $ws.Range('E19')
$cell.NumberFormat = '@'
$cell.Value = '  -0.71%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D20')
$cell.NumberFormat = '@'
$cell.Value = '352.02'
$cell.Style = 'Normal'

$cell = $ws.Range('E20')
$cell.NumberFormat = '@'
$cell.Value = '  +0.28%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D21')
$cell.NumberFormat = '@'
$cell.Value = '4.03'
$cell.Style = 'Normal'

$cell = $ws.Range('E21')
$cell.NumberFormat = '@'
$cell.Value = '  -0.45%  '
$cell.Style = 'Normal'

$cell = $ws.Range('E22')
$cell.NumberFormat = '@'
$cell.Value = '  -0.03%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D23')
$cell.NumberFormat = '@'
$cell.Value = '70.52'
$cell.Style = 'Normal'

$cell = $ws.Range('E23')
$cell.NumberFormat = '@'
$cell.Value = '  +2.37%  '
$cell.Style = 'Normal'

$cell = $ws.Range('E24')
$cell.NumberFormat = '@'
$cell.Value = '  -0.45%  '
$cell.Style = 'Normal'

$cell = $ws.Range('E25')
$cell.NumberFormat = '@'
$cell.Value = '  -0.98%  '
$cell.Style = 'Normal'

$cell = $ws.Range('E26')
$cell.NumberFormat = '@'
$cell.Value = '  +1.05%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D27')
$cell.NumberFormat = '@'
$cell.Value = '2.615.50'
$cell.Style = 'Normal'

$cell = $ws.Range('E27')
$cell.NumberFormat = '@'
$cell.Value = '  +0.74%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D28')
$cell.NumberFormat = '@'
$cell.Value = '0.998'
$cell.Style = 'Normal'

$cell = $ws.Range('E28')
$cell.NumberFormat = '@'
$cell.Value = '  -0.55%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D30')
$cell.NumberFormat = '@'
$cell.Value = '508.61'
$cell.Style = 'Normal'

$cell = $ws.Range('E30')
$cell.NumberFormat = '@'
$cell.Value = '  -0.09%  '
$cell.Style = 'Normal'

$cell = $ws.Range('E31')
$cell.NumberFormat = '@'
$cell.Value = '  +1.79%  '
$cell.Style = 'Normal'

$cell = $ws.Range('E32')
$cell.NumberFormat = '@'
$cell.Value = '  +2.00%  '
$cell.Style = 'Normal'

$cell = $ws.Range('E33')
$cell.NumberFormat = '@'
$cell.Value = '  +0.52%  '
$cell.Style = 'Normal'

$cell = $ws.Range('E34')
$cell.NumberFormat = '@'
$cell.Value = '  -0.01%  '
$cell.Style = 'Normal'

$cell = $ws.Range('E35')
$cell.NumberFormat = '@'
$cell.Value = '  +5.76%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D36')
$cell.NumberFormat = '@'
$cell.Value = '161.04'
$cell.Style = 'Normal'

$cell = $ws.Range('E36')
$cell.NumberFormat = '@'
$cell.Value = '  +1.01%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D37')
$cell.NumberFormat = '@'
$cell.Value = '18.69'
$cell.Style = 'Normal'

$cell = $ws.Range('E37')
$cell.NumberFormat = '@'
$cell.Value = '  +0.17%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D38')
$cell.NumberFormat = '@'
$cell.Value = '18.33'
$cell.Style = 'Normal'

$cell = $ws.Range('E38')
$cell.NumberFormat = '@'
$cell.Value = '  +0.06%  '
$cell.Style = 'Normal'

$cell = $ws.Range('E39')
$cell.NumberFormat = '@'
$cell.Value = '  +0.67%  '
$cell.Style = 'Normal'

$cell = $ws.Range('E40')
$cell.NumberFormat = '@'
$cell.Value = '  +0.01%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D41')
$cell.NumberFormat = '@'
$cell.Value = '1.71'
$cell.Style = 'Normal'

$cell = $ws.Range('E41')
$cell.NumberFormat = '@'
$cell.Value = '  +1.34%  '
$cell.Style = 'Normal'

$cell = $ws.Range('E42')
$cell.NumberFormat = '@'
$cell.Value = '  +0.85%  '
$cell.Style = 'Normal'

$cell = $ws.Range('E43')
$cell.NumberFormat = '@'
$cell.Value = '  +1.31%  '
$cell.Style = 'Normal'

$cell = $ws.Range('E44')
$cell.NumberFormat = '@'
$cell.Value = '  +2.64%  '
$cell.Style = 'Normal'

$cell = $ws.Range('D45')
$cell.NumberFormat = '@'
$cell.Value = '143.64'
$cell.Style = 'Normal'

$cell = $ws.Range('E45')
$cell.NumberFormat = '@'
$cell.Value = '  +1.98%  '
$cell.Style = 'Normal'

$cell = $ws.Range('E46')
$cell.NumberFormat = '@'
$cell.Value = '  +1.86%  '
$cell.Style = 'Normal'

$cell = $ws.Range('B47')
$cell.NumberFormat = '@'
$cell.Value = 'ARBITRUM'
$cell.Style = 'Normal'

$cell = $ws.Range('C47')
$cell.NumberFormat = '@'
$cell.Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$cell.Style = 'Normal'

$cell = $ws.Range('D47')
$cell.NumberFormat = '@'
$cell.Value = '0.513'
$cell.Style = 'Normal'

$cell = $ws.Range('E47')
$cell.NumberFormat = '@'
$cell.Value = '  +0.09%  '
$cell.Style = 'Normal'

$cell = $ws.Range('B48')
$cell.NumberFormat = '@'
$cell.Value = 'Cronos'
$cell.Style = 'Normal'

$cell = $ws.Range('C48')
$cell.NumberFormat = '@'
$cell.Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$cell.Style = 'Normal'

$cell = $ws.Range('D48')
$cell.NumberFormat = '@'
$cell.Value = '0.0745'
$cell.Style = 'Normal'

$cell = $ws.Range('E48')
$cell.NumberFormat = '@'
$cell.Value = '  +2.08%  '
$cell.Style = 'Normal'

$cell = $ws.Range('B49')
$cell.NumberFormat = '@'
$cell.Value = 'Optimism'
$cell.Style = 'Normal'

$cell = $ws.Range('C49')
$cell.NumberFormat = '@'
$cell.Value = 'https://coinranking.com/coin/n1p-s_gm1+optimism-op'
$cell.Style = 'Normal'

$cell = $ws.Range('D49')
$cell.NumberFormat = '@'
$cell.Value = '1.58'
$cell.Style = 'Normal'

$cell = $ws.Range('E49')
$cell.NumberFormat = '@'
$cell.Value = '  -0.22%  '
$cell.Style = 'Normal'

$cell = $ws.Range('B50')
$cell.NumberFormat = '@'
$cell.Value = 'Mantle'
$cell.Style = 'Normal'

$cell = $ws.Range('C50')
$cell.NumberFormat = '@'
$cell.Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$cell.Style = 'Normal'

$cell = $ws.Range('D50')
$cell.NumberFormat = '@'
$cell.Value = '0.587'
$cell.Style = 'Normal'

$cell = $ws.Range('E50')
$cell.NumberFormat = '@'
$cell.Value = '  +1.01%  '
$cell.Style = 'Normal'

$cell = $ws.Range('B51')
$cell.NumberFormat = '@'
$cell.Value = 'BitgetToken'
$cell.Style = 'Normal'

$cell = $ws.Range('C51')
$cell.NumberFormat = '@'
$cell.Value = 'https://coinranking.com/coin/q7gMmMdLb+bitgettoken-bgb'
$cell.Style = 'Normal'

$cell = $ws.Range('D51')
$cell.NumberFormat = '@'
$cell.Value = '1.19'
$cell.Style = 'Normal'

$cell = $ws.Range('E51')
$cell.NumberFormat = '@'
$cell.Value = '  +1.82%  '
$cell.Style = 'Normal'

